# Apply the cmd_options.xlsx matrix-simplification edit:
#  - Row 14/15 style + value tweaks (new unpack_argument/mapped_key truth table)
#  - Remove rows 16-21 (six now-obsolete combinations)
#  - Rows 27-30 shift up to 21-24 automatically once 16-21 are removed
#  - Move the active selection to B19
#  - Tidy the absPath note (best effort; not exposed via the Excel object model)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats - used to copy just the cell style (fill/alignment) between cells
# without disturbing the destination's existing value.
$xlPasteFormats = -4122

# --- Row 14 -----------------------------------------------------------
# C14: red -> green
$ws.Range("B2").Copy()
$ws.Range("C14").PasteSpecial($xlPasteFormats)
# E14: green -> "none" fill (style 4)
$ws.Range("J27").Copy()
$ws.Range("E14").PasteSpecial($xlPasteFormats)
# K14: gains the checkmark value (reuses the existing shared string)
$ws.Range("K14").Value = "✅"

# --- Row 15 -----------------------------------------------------------
# C15: red -> green
$ws.Range("B2").Copy()
$ws.Range("C15").PasteSpecial($xlPasteFormats)
# E15: red -> "none" fill (style 4)
$ws.Range("J27").Copy()
$ws.Range("E15").PasteSpecial($xlPasteFormats)
# F15: green -> red
$ws.Range("F2").Copy()
$ws.Range("F15").PasteSpecial($xlPasteFormats)
# G15: cell removed entirely
$ws.Range("G15").Clear()
# P15: loses its checkmark value but keeps its style
$ws.Range("P15").ClearContents()

$excel.CutCopyMode = $false

# --- Remove now-obsolete rows 16-21 ------------------------------------
# This also shifts the old rows 27-30 up to become 21-24.
$ws.Range("16:21").Delete()

# --- Update the active selection ---------------------------------------
$ws.Range("B19").Select()
